# "use new word template" (#4035)
#
# The refreshed blank-document template renames the four built-in default
# styles' internal style IDs (Word keeps the display name / w:name the
# same, only the short w:styleId token changes):
#
#   Normal               -> a
#   DefaultParagraphFont  -> a0
#   TableNormal           -> a1
#   NoList                -> a2
#
# Word's object model does not expose a settable Style.StyleId, so the
# rename is performed by re-creating each style under its new id (via
# Styles.Add, whose id is derived from the name with spaces stripped)
# and then restoring the original display name / flags, before removing
# the old entries. Styles.Add appends, so re-adding in original order
# and deleting the stale originals afterwards keeps the stylesheet in
# its original relative order.

$d = $word.ActiveDocument

$normal = $d.Styles.Add("a", 1)
$normal.NameLocal = "Normal"
$normal.QuickStyle = $true

$defaultParagraphFont = $d.Styles.Add("a0", 2)
$defaultParagraphFont.NameLocal = "Default Paragraph Font"
$defaultParagraphFont.Priority = 1
$defaultParagraphFont.UnhideWhenUsed = $true

$tableNormal = $d.Styles.Add("a1", 3)
$tableNormal.NameLocal = "Normal Table"
$tableNormal.Priority = 99
$tableNormal.UnhideWhenUsed = $true

$noList = $d.Styles.Add("a2", 4)
$noList.NameLocal = "No List"
$noList.Priority = 99
$noList.UnhideWhenUsed = $true

# Drop the stale built-ins now that their replacements are in place.
$d.Styles.Item("Normal").Delete()
$d.Styles.Item("DefaultParagraphFont").Delete()
$d.Styles.Item("TableNormal").Delete()
$d.Styles.Item("NoList").Delete()
